# Add two new columns (I: "I0", J: "IF") to the sheet, mirroring the
# existing header style used by the other header cells (e.g. H1), and
# fill in the corresponding numeric values for data rows 2-75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of the existing "IP" header cell (H1) onto
# the two new header cells so they pick up the same bold/centered/bordered
# style already used by the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows (2-75) -------------------------------------------------
$data = @(
    @{ Row = 2; I = 9; J = 9 },
    @{ Row = 3; I = 7; J = 7 },
    @{ Row = 4; I = 9; J = 9 },
    @{ Row = 5; I = 8; J = 9 },
    @{ Row = 6; I = 8; J = 9 },
    @{ Row = 7; I = 9; J = 9 },
    @{ Row = 8; I = 9; J = 9 },
    @{ Row = 9; I = 10; J = 10 },
    @{ Row = 10; I = 9; J = 9 },
    @{ Row = 11; I = 8; J = 9 },
    @{ Row = 12; I = 8; J = 8 },
    @{ Row = 13; I = 9; J = 9 },
    @{ Row = 14; I = 9; J = 9 },
    @{ Row = 15; I = 9; J = 9 },
    @{ Row = 16; I = 9; J = 9 },
    @{ Row = 17; I = 9; J = 9 },
    @{ Row = 18; I = 9; J = 9 },
    @{ Row = 19; I = 9; J = 9 },
    @{ Row = 20; I = 7; J = 7 },
    @{ Row = 21; I = 7; J = 7 },
    @{ Row = 22; I = 7; J = 7 },
    @{ Row = 23; I = 8; J = 8 },
    @{ Row = 24; I = 9; J = 9 },
    @{ Row = 25; I = 9; J = 9 },
    @{ Row = 26; I = 9; J = 9 },
    @{ Row = 27; I = 10; J = 10 },
    @{ Row = 28; I = 9; J = 9 },
    @{ Row = 29; I = 9; J = 9 },
    @{ Row = 30; I = 9; J = 9 },
    @{ Row = 31; I = 9; J = 9 },
    @{ Row = 32; I = 9; J = 9 },
    @{ Row = 33; I = 9; J = 9 },
    @{ Row = 34; I = 9; J = 10 },
    @{ Row = 35; I = 9; J = 9 },
    @{ Row = 36; I = 9; J = 9 },
    @{ Row = 37; I = 8; J = 9 },
    @{ Row = 38; I = 9; J = 9 },
    @{ Row = 39; I = 8; J = 8 },
    @{ Row = 40; I = 9; J = 9 },
    @{ Row = 41; I = 10; J = 10 },
    @{ Row = 42; I = 7; J = 8 },
    @{ Row = 43; I = 9; J = 9 },
    @{ Row = 44; I = 9; J = 9 },
    @{ Row = 45; I = 9; J = 9 },
    @{ Row = 46; I = 9; J = 9 },
    @{ Row = 47; I = 10; J = 10 },
    @{ Row = 48; I = 9; J = 9 },
    @{ Row = 49; I = 8; J = 9 },
    @{ Row = 50; I = 8; J = 9 },
    @{ Row = 51; I = 9; J = 9 },
    @{ Row = 52; I = 9; J = 9 },
    @{ Row = 53; I = 9; J = 9 },
    @{ Row = 54; I = 9; J = 9 },
    @{ Row = 55; I = 9; J = 9 },
    @{ Row = 56; I = 9; J = 9 },
    @{ Row = 57; I = 9; J = 9 },
    @{ Row = 58; I = 9; J = 9 },
    @{ Row = 59; I = 8; J = 9 },
    @{ Row = 60; I = 8; J = 8 },
    @{ Row = 61; I = 9; J = 9 },
    @{ Row = 62; I = 9; J = 9 },
    @{ Row = 63; I = 9; J = 9 },
    @{ Row = 64; I = 8; J = 8 },
    @{ Row = 65; I = 7; J = 7 },
    @{ Row = 66; I = 9; J = 9 },
    @{ Row = 67; I = 9; J = 9 },
    @{ Row = 68; I = 7; J = 7 },
    @{ Row = 69; I = 8; J = 8 },
    @{ Row = 70; I = 8; J = 8 },
    @{ Row = 71; I = 6; J = 6 },
    @{ Row = 72; I = 6; J = 6 },
    @{ Row = 73; I = 6; J = 6 },
    @{ Row = 74; I = 4; J = 4 },
    @{ Row = 75; I = 3; J = 3 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I   # column I
    $ws.Cells.Item($r, 10).Value = $entry.J  # column J
}

Write-Host "Added columns I (I0) and J (IF) for rows 1-75"
